$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C12").Value = 8508
$ws.Range("C13:C25").Value = 8448
$ws.Range("C26:C54").Value = 8415
$ws.Range("C55:C58").Value = 8184
$ws.Range("C59:C76").Value = 7594
$ws.Range("C77:C123").Value = 7569
$ws.Range("C124:C137").Value = 7295
$ws.Range("C141:C169").Value = 7293
